# Append: 2025-09-29 12:48 JST
# Update the "取得日時" (retrieved datetime) timestamps in column A of the
# "ランサーズ" (Lancers) sheet for the currently listed rows (2-9) from the
# old scrape time to the new scrape time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-09-29 12:36:28"
$newTimestamp = "2025-09-29 12:48:51"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
